# Apply the diff: update the date line and the 25 division-problem answer cells.
$d = $word.ActiveDocument

# 1. Update the date paragraph at the top of the document.
$d.Content.Find.Execute("2024-11-15 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-11-16 Saturday", 2)

# 2. Update the answer text in each populated table cell (row/col are 1-indexed).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "61÷6=10, 1"
$t.Cell(1, 2).Range.Text = "81÷6=13, 3"
$t.Cell(1, 3).Range.Text = "10÷7=1, 3"
$t.Cell(1, 4).Range.Text = "60÷2=30, 0"
$t.Cell(1, 5).Range.Text = "57÷4=14, 1"

$t.Cell(5, 1).Range.Text = "40÷6=6, 4"
$t.Cell(5, 2).Range.Text = "75÷8=9, 3"
$t.Cell(5, 3).Range.Text = "61÷5=12, 1"
$t.Cell(5, 4).Range.Text = "54÷9=6, 0"
$t.Cell(5, 5).Range.Text = "72÷5=14, 2"

$t.Cell(9, 1).Range.Text = "56÷7=8, 0"
$t.Cell(9, 2).Range.Text = "83÷5=16, 3"
$t.Cell(9, 3).Range.Text = "29÷2=14, 1"
$t.Cell(9, 4).Range.Text = "86÷7=12, 2"
$t.Cell(9, 5).Range.Text = "45÷2=22, 1"

$t.Cell(13, 1).Range.Text = "24÷4=6, 0"
$t.Cell(13, 2).Range.Text = "57÷9=6, 3"
$t.Cell(13, 3).Range.Text = "33÷8=4, 1"
$t.Cell(13, 4).Range.Text = "66÷6=11, 0"
$t.Cell(13, 5).Range.Text = "85÷9=9, 4"

$t.Cell(17, 1).Range.Text = "87÷6=14, 3"
$t.Cell(17, 2).Range.Text = "70÷8=8, 6"
$t.Cell(17, 3).Range.Text = "52÷5=10, 2"
$t.Cell(17, 4).Range.Text = "69÷2=34, 1"
$t.Cell(17, 5).Range.Text = "97÷4=24, 1"

Write-Host "Done applying edits."
